$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "MuSCs"
$ws.Range("G2").Value2 = 3.855689333333333
$ws.Range("I2").Value2 = 0.1513312545414155
$ws.Range("J2").Value2 = 0.1513312545414155
$ws.Range("M2").Value2 = 0.08675100000000001
$ws.Range("N2").Value2 = 0.260253
$ws.Range("O2").Value2 = 0.2009041159973413
$ws.Range("P2").Value2 = 0.2009041159973413
$ws.Range("Q2").Value2 = 0.334484905356
$ws.Range("R2").Value2 = 3.010364148204
$ws.Range("S2").Value2 = 0.03040307191641172
$ws.Range("T2").Value2 = 0.03040307191641173

$ws.Range("D3").Value2 = "Resolving-Mac"
$ws.Range("G3").Value2 = 3.855689333333333
$ws.Range("I3").Value2 = 0.1513312545414155
$ws.Range("J3").Value2 = 0.1513312545414155
$ws.Range("K3").Value2 = 3
$ws.Range("L3").Value2 = 1
$ws.Range("M3").Value2 = 0.3450520000000001
$ws.Range("N3").Value2 = 1.035156
$ws.Range("O3").Value2 = 0.7990958840026586
$ws.Range("P3").Value2 = 0.7990958840026586
$ws.Range("Q3").Value2 = 1.330413315845334
$ws.Range("R3").Value2 = 11.973719842608
$ws.Range("S3").Value2 = 0.1209281826250037
$ws.Range("T3").Value2 = 0.1209281826250038

$ws.Range("A4").Value2 = "MuSCs"
$ws.Range("D4").Value2 = "MuSCs"
$ws.Range("G4").Value2 = 5.360192666666666
$ws.Range("H4").Value2 = 16.080578
$ws.Range("I4").Value2 = 0.2103812342497758
$ws.Range("J4").Value2 = 0.2103812342497758
$ws.Range("K4").Value2 = 1
$ws.Range("L4").Value2 = 0.3333333333333333
$ws.Range("M4").Value2 = 0.08675100000000001
$ws.Range("N4").Value2 = 0.260253
$ws.Range("O4").Value2 = 0.2009041159973413
$ws.Range("P4").Value2 = 0.2009041159973413
$ws.Range("Q4").Value2 = 0.465002074026
$ws.Range("R4").Value2 = 4.185018666234
$ws.Range("S4").Value2 = 0.0422664558893808
$ws.Range("T4").Value2 = 0.04226645588938081

$ws.Range("D5").Value2 = "Resolving-Mac"
$ws.Range("G5").Value2 = 5.360192666666666
$ws.Range("H5").Value2 = 16.080578
$ws.Range("I5").Value2 = 0.2103812342497758
$ws.Range("J5").Value2 = 0.2103812342497758
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.3450520000000001
$ws.Range("N5").Value2 = 1.035156
$ws.Range("O5").Value2 = 0.7990958840026586
$ws.Range("P5").Value2 = 0.7990958840026586
$ws.Range("Q5").Value2 = 1.849545200018667
$ws.Range("R5").Value2 = 16.645906800168
$ws.Range("S5").Value2 = 0.168114778360395
$ws.Range("T5").Value2 = 0.168114778360395

$ws.Range("A6").Value2 = "Resolving-Mac"
$ws.Range("G6").Value2 = 16.26259133333333
$ws.Range("H6").Value2 = 48.787774
$ws.Range("I6").Value2 = 0.6382875112088087
$ws.Range("J6").Value2 = 0.6382875112088088
$ws.Range("K6").Value2 = 1
$ws.Range("L6").Value2 = 0.3333333333333333
$ws.Range("M6").Value2 = 0.08675100000000001
$ws.Range("N6").Value2 = 0.260253
$ws.Range("O6").Value2 = 0.2009041159973413
$ws.Range("P6").Value2 = 0.2009041159973413
$ws.Range("Q6").Value2 = 1.410796060758
$ws.Range("R6").Value2 = 12.697164546822
$ws.Range("S6").Value2 = 0.1282345881915488
$ws.Range("T6").Value2 = 0.1282345881915488

$ws.Range("A7").Value2 = "Resolving-Mac"
$ws.Range("G7").Value2 = 16.26259133333333
$ws.Range("H7").Value2 = 48.787774
$ws.Range("I7").Value2 = 0.6382875112088087
$ws.Range("J7").Value2 = 0.6382875112088088
$ws.Range("M7").Value2 = 0.3450520000000001
$ws.Range("N7").Value2 = 1.035156
$ws.Range("O7").Value2 = 0.7990958840026586
$ws.Range("P7").Value2 = 0.7990958840026586
$ws.Range("Q7").Value2 = 5.611439664749335
$ws.Range("R7").Value2 = 50.50295698274401
$ws.Range("S7").Value2 = 0.5100529230172598
$ws.Range("T7").Value2 = 0.5100529230172599

$ws.Rows("8:10").Delete()
